$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.694.21"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "2.896.73"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'354.24"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'108.85"
$ws.Range("E6").Value = "  -3.24%  "
$ws.Range("D7").Value = "'0.561"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("D10").Value = "'38.71"
$ws.Range("E10").Value = "  -3.51%  "
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "'19.38"
$ws.Range("E13").Value = "  -3.06%  "
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "3.358.42"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "2.903.69"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "'0.971"
$ws.Range("D18").Value = "51.629.86"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("D20").Value = "'7.49"
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("D21").Value = "'13.74"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("D24").Value = "'267.10"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").Value = "'2.78"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("E26").Value = "  +9.11%  "
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("E28").Value = "  +16.00%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +8.29%  "
$ws.Range("D31").Value = "'10.43"
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("D32").Value = "'37.22"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("D33").Value = "'2.20"
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("D34").Value = "'6.09"
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("D35").Value = "'51.99"
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("E36").Value = "  -3.06%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "'3.16"
$ws.Range("E38").Value = "  -4.77%  "
$ws.Range("D39").Value = "'18.10"
$ws.Range("E39").Value = "  -3.86%  "
$ws.Range("D40").Value = "'1.99"
$ws.Range("E40").Value = "  -4.10%  "
$ws.Range("E41").Value = "  -7.70%  "
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("D43").Value = "'22.56"
$ws.Range("E43").Value = "  -4.18%  "
$ws.Range("D44").Value = "'118.62"
$ws.Range("E44").Value = "  -2.57%  "
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("E46").Value = "  -5.86%  "
$ws.Range("D47").Value = "'3.42"
$ws.Range("E47").Value = "  -4.27%  "
$ws.Range("D48").Value = "2.118.76"
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("E49").Value = "  -6.31%  "
$ws.Range("D50").Value = "'0.0336"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("E51").Value = "  -0.69%  "

Write-Host "Applied crypto updates"
